$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new row above row 5 ("신한스팩13호" row), shifting everything down.
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the new SPAC entry.
$ws.Range("A5").Value = "하나스팩33호"
$ws.Range("B5").Value = "2024.04.08~04.09"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 7000
$ws.Range("F5").Value = "하나증권"

# 하나스팩32호's IPO price was confirmed (was "-", now 2000). After the
# insert above, that row (originally row 8) is now row 9.
$ws.Range("D9").Value = "2000"

# Drop the now-duplicated last row (old row 21 shifted to row 22) to keep the
# table at 20 data rows (A1:F21).
$ws.Rows.Item(22).Delete()
